# Update cryptos list (prices + volume%) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold text like "30.621.37" or "0.06686" -- force
# Text format before assigning so Excel does not reinterpret them as numbers
# (which would normalize "0.06600" -> 0.066, drop thousands-style dots, etc.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.526.21"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.884.34"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.71%  "

$ws.Range("D5").Value = "235.83"
$ws.Range("E5").Value = "  -3.67%  "

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").Value = "0.4847"
$ws.Range("E7").Value = "  -1.63%  "

$ws.Range("D8").Value = "0.2890"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").Value = "0.06600"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "1.889.96"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").Value = "16.71"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "0.07193"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").Value = "88.36"
$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").Value = "4.964"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("D15").Value = "0.6620"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("D16").Value = "30.527.57"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "0.000007818"
$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "2.137.83"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").Value = "4.742"
$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("D23").Value = "186.06"
$ws.Range("E23").Value = "  +25.32%  "

$ws.Range("D24").Value = "5.980"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "9.216"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").Value = "155.04"
$ws.Range("E26").Value = "  +2.04%  "

$ws.Range("D27").Value = "18.45"
$ws.Range("E27").Value = "  +7.84%  "

$ws.Range("D28").Value = "1.844"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "4.219"
$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("D31").Value = "0.08955"
$ws.Range("E31").Value = "  +1.52%  "

$ws.Range("D32").Value = "3.895"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").Value = "0.05225"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").Value = "0.7259"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.749"
$ws.Range("E35").Value = "  +3.10%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.076"
$ws.Range("E36").Value = "  -4.22%  "

$ws.Range("D37").Value = "0.01808"
$ws.Range("E37").Value = "  -2.01%  "

$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").Value = "0.9200"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("D40").Value = "2.052"
$ws.Range("E40").Value = "  -6.06%  "

$ws.Range("D41").Value = "0.4306"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "103.72"
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").Value = "5.593"
$ws.Range("E44").Value = "  -4.13%  "

$ws.Range("D45").Value = "0.1327"
$ws.Range("E45").Value = "  +1.69%  "

$ws.Range("D46").Value = "7.325"
$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("D47").Value = "0.05845"
$ws.Range("E47").Value = "  +0.96%  "

$ws.Range("D48").Value = "8.715"
$ws.Range("E48").Value = "  +4.51%  "

$ws.Range("D49").Value = "33.23"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("D50").Value = "0.3878"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("D51").Value = "1.396"
$ws.Range("E51").Value = "  +3.08%  "
